$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.18345320224762
$ws.Range("B1").Value = 2.298865795135498
$ws.Range("C1").Value = 4.549642086029053
$ws.Range("D1").Value = 3.455256462097168
$ws.Range("E1").Value = 1.205423831939697
